$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MSG: None

MSG: The rights to both movies have been successfully acquired.
"
$ws.Range("D2").Value = "both_movies, "
$ws.Range("C3").Value = "MSG: None

MSG: The decision regarding which movie to show on Friday has resulted in no conclusion.
"
$ws.Range("D3").Value = "no_decision, "
$ws.Range("C4").Value = "MSG: None

MSG: The rights to both movies have been successfully acquired.
"
$ws.Range("C5").Value = "MSG: None

MSG: The decision concluded with no movie selected for Friday.
"
$ws.Range("C6").Value = "MSG: None

MSG: The decision has been recorded, reflecting that no movie was selected for Friday.
"
$ws.Range("C7").Value = "MSG: None

MSG: The decision has been recorded as no movie selected.
"
$ws.Range("C8").Value = "MSG: None

MSG: The decision has been recorded: `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("C9").Value = "MSG: None

MSG: The decision has been recorded as no choice of a movie was made.
"
$ws.Range("C10").Value = "MSG: None

MSG: The decision has been made to acquire the rights to the movie `"Barbie.`"
"
$ws.Range("C11").Value = "MSG: None

MSG: The decision has been recorded, and `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("C12").Value = "MSG: None

MSG: The decision has been recorded, indicating that no agreement was reached regarding the movie to be shown on Friday.
"
$ws.Range("C13").Value = "MSG: None

MSG: The committee did not reach a decision on which movie to show on Friday.
"
$ws.Range("C14").Value = "MSG: None

MSG: The movie `"Barbie`" has been successfully selected for acquisition.
"
$ws.Range("C15").Value = "MSG: None

MSG: The decision-making committee did not reach a conclusion regarding which movie to show on Friday, resulting in no decision being made.
"
$ws.Range("C16").Value = "MSG: None

MSG: I have recorded the decision as no decision was made regarding the movie for Friday.
"
$ws.Range("C17").Value = "MSG: None

MSG: The decision has been made to acquire rights for `"Barbie`" to be shown on Friday.
"
$ws.Range("C18").Value = "MSG: None

MSG: The decision has been recorded as no choice for Friday's movie could be made.
"
$ws.Range("C19").Value = "MSG: None

MSG: The function has been called, indicating that no decision about Friday's movie has been made.
"
$ws.Range("C20").Value = "MSG: None

MSG: The decision on which movie to show on Friday has not been made.
"
$ws.Range("C21").Value = "MSG: None

MSG: The decision regarding the movie selection for Friday has ended without a consensus, resulting in no decision being made.
"
$ws.Range("C22").Value = "MSG: None

MSG: The decision has been recorded as no movie selected for Friday.
"
$ws.Range("C23").Value = "MSG: None

MSG: The committee did not reach a decision on the movie to be shown on Friday.
"
$ws.Range("C24").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie`" as the movie to be shown on Friday.
"
$ws.Range("C25").Value = "MSG: None

MSG: The decision has been recorded, indicating that no movie was selected during the discussion.
"
$ws.Range("C26").Value = "MSG: None

MSG: The decision to acquire a movie was not finalized, and thus no action will be taken regarding the movie selection at this time.
"
$ws.Range("D26").Value = "no_decision, "
$ws.Range("C27").Value = "MSG: None

MSG: The decision has been recorded, and no movie was selected for Friday.
"
$ws.Range("C28").Value = "MSG: None

MSG: The decision regarding the movie to show on Friday has not reached a conclusion.
"
$ws.Range("C29").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie.`"
"
$ws.Range("C30").Value = "MSG: None

MSG: The rights to `"Barbie`" have been successfully acquired for the Friday screening.
"
$ws.Range("C31").Value = "MSG: None

MSG: The decision has been recorded, indicating that no agreement was reached regarding the movie selection for Friday.
"
$ws.Range("C32").Value = "MSG: None

MSG: The decision has been recorded as `"no_decision.`" There was no agreement reached on a movie to be shown on Friday.
"
$ws.Range("C33").Value = "MSG: None

MSG: The movie `"Barbie`" has been selected for Friday's showing.
"
$ws.Range("C34").Value = "MSG: None

MSG: The decision about Friday's movie cannot be made at this time.
"
$ws.Range("C35").Value = "MSG: None

MSG: The decision has been recorded as no decision.
"
$ws.Range("C36").Value = "MSG: None

MSG: The decision has been recorded, and `"Barbie`" has been selected for acquisition.
"
$ws.Range("C37").Value = "MSG: None

MSG: No decision was made regarding which movie to show on Friday.
"
$ws.Range("C38").Value = "MSG: None

MSG: The decision was made to not acquire a movie for Friday, as no consensus was reached.
"
$ws.Range("C39").Value = "MSG: None

MSG: The rights to both movies have been acquired successfully for the showing on Friday.
"
$ws.Range("C40").Value = "MSG: None

MSG: The decision to select `"Barbie`" has been recorded.
"
$ws.Range("C41").Value = "MSG: None

MSG: The decision has been recorded with no movie selected for Friday.
"
$ws.Range("C42").Value = "MSG: None

MSG: The decision has been recorded as no decision was reached regarding the movie for Friday.
"
$ws.Range("C43").Value = "MSG: None

MSG: The decision process has concluded without a definitive choice for Friday's movie.
"
$ws.Range("C44").Value = "MSG: None

MSG: The decision has been recorded as no agreement was reached regarding the movie to be shown on Friday.
"
$ws.Range("C45").Value = "MSG: None

MSG: The rights to both movies have been successfully acquired.
"
$ws.Range("C46").Value = "MSG: None

MSG: The decision has been recorded, marking that no movie was selected for the assembly.
"
$ws.Range("C47").Value = "MSG: None

MSG: The decision has been recorded, confirming that `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("C48").Value = "MSG: None

MSG: The committee did not reach a decision about which movie to show on Friday.
"
$ws.Range("C49").Value = "MSG: None

MSG: The decision has been recorded as `"Barbie`" being selected for the movie on Friday.
"
$ws.Range("C50").Value = "MSG: None

MSG: The decision has been recorded as no movie selected for Friday, as the committee could not reach a consensus.
"
$ws.Range("C51").Value = "MSG: None

MSG: The decision to acquire the rights for `"Barbie`" has been successfully recorded.
"
$ws.Range("C52").Value = "MSG: None

MSG: The function for no decision has been called, indicating that no choice of a movie is possible without further discussion.
"
$ws.Range("C53").Value = "MSG: None

MSG: The decision regarding the movie to be shown on Friday could not be made, so no actions will be taken.
"
$ws.Range("C54").Value = "MSG: None

MSG: The decision has been recorded, and there is no selected movie for Friday.
"
$ws.Range("C55").Value = "MSG: None

MSG: The decision has been made to acquire the rights for `"Barbie`" to be shown on Friday.
"
$ws.Range("D55").Value = "Barbie_was_selected, "
$ws.Range("C56").Value = "MSG: None

MSG: The decision has been recorded, indicating that no choice for Friday's movie was made.
"
$ws.Range("C57").Value = "MSG: None

MSG: The decision has been recorded. `"Barbie`" will be the movie shown on Friday.
"
$ws.Range("C58").Value = "MSG: None

MSG: No decision was made regarding the movie to be shown on Friday.
"
$ws.Range("C59").Value = "MSG: None

MSG: The decision has been recorded as no decision was made regarding which movie to show on Friday.
"
$ws.Range("C60").Value = "MSG: None

MSG: The decision has been recorded with no agreement on which movie to acquire for Friday.
"
$ws.Range("C61").Value = "MSG: None

MSG: The committee did not reach a decision regarding which movie to show on Friday.
"
$ws.Range("C62").Value = "MSG: None

MSG: The decision has been recorded, and the rights to `"Barbie`" will be acquired for the upcoming showing.
"
$ws.Range("C63").Value = "MSG: None

MSG: The decision has been recorded to acquire the rights for `"Barbie.`"
"
$ws.Range("C64").Value = "MSG: None

MSG: I have successfully recorded the decision to acquire the rights for both movies.
"
$ws.Range("C65").Value = "MSG: None

MSG: The decision-making process ended without a selection for Friday's movie, leading to no decision being made.
"
$ws.Range("C66").Value = "MSG: None

MSG: The decision regarding which movie to show on Friday resulted in no selection being made.
"
$ws.Range("C67").Value = "MSG: None

MSG: The decision has been recorded as `"no decision.`"
"
$ws.Range("C68").Value = "MSG: None

MSG: The decision is that no movie was selected.
"
